# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 05:52"

# 2) Refresh Australia's stats (row 33: rank 37)
$ws.Range("B33").Value = 6394
$ws.Range("C33").Value = 35
$ws.Range("E33").Value = 2839

# 3) Insert Mongolia into its correctly sorted position (just above Sudan,
#    currently row 164) with its refreshed figures, pushing Sudan and the
#    following countries down by one row.
$ws.Rows("164:164").Insert()
$ws.Range("A164").Value = "Mongolia"
$ws.Range("B164").Value = 30
$ws.Range("C164").Value = 13
$ws.Range("D164").Value = 11
$ws.Range("E164").Value = 19
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 0

# Remove the old Mongolia row, which shifted down to row 179 after the
# insert above, restoring the table to its original length.
$ws.Rows("179:179").Delete()
